# Replace the remaining "bootstrap" numeric placeholder values with the
# new "material ui" values in the insulation-resistance measurement table.
#
# Each old value occurs exactly once in the document, so straightforward
# Find/Replace (matching the whole word, i.e. the complete numeric token)
# is sufficient and avoids any ordering collisions between old/new values
# that coincidentally overlap (e.g. 0.314 appears both as an old and a
# new value in different cells).

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "0.059"; New = "0.161" },
    @{ Old = "0.702"; New = "0.435" },
    @{ Old = "0.438"; New = "0.695" },
    @{ Old = "0.314"; New = "0.395" },
    @{ Old = "0.606"; New = "0.492" },
    @{ Old = "0.162"; New = "0.347" },
    @{ Old = "0.192"; New = "0.734" },
    @{ Old = "0.52";  New = "0.216" },
    @{ Old = "0.846"; New = "0.468" },
    @{ Old = "0.968"; New = "0.697" },
    @{ Old = "0.141"; New = "0.314" },
    @{ Old = "0.233"; New = "0.942" }
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.Old, $true, $true, $false, $false, $false, $true, 1, $false, $pair.New, 2)
}
